$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new trailing header columns on the frozen title row.
$ws.Cells.Item(1, 14).Value = "評估淨值"
$ws.Cells.Item(1, 15).Value = "貸放成數"

# Match the saved selection/zoom state recorded in the sheet view.
$ws.Range("N1:O1").Select()
$excel.ActiveWindow.Zoom = 55
